# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Update "Bad Drivers" table values
$ws.Range("D3").Value = 98.90000000000001
$ws.Range("C4").Value = 355
$ws.Range("D4").Value = 98.90000000000001
$ws.Range("C5").Value = 363
